$d = $word.ActiveDocument
$d.Content.Find.Execute("Nowy plik", $false, $false, $false, $false, $false, $true, 1, $false, "Zmiana zmiana zmiana", 2)
